$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-33 change from serial date 45180 to 45181
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value = 45181
}
